$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Range("B29").Value = 6865295
$ws.Range("C29").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D29").Value = 45172.61458333334
$ws.Range("E29").Value = 'FK Tuzla City'
$ws.Range("F29").Value = 'NK Igman Konjic'
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 1
$ws.Range("K29").Value = 'H'
$ws.Range("L29").Value = 1.8
$ws.Range("M29").Value = 3.4
$ws.Range("N29").Value = 3.8
$ws.Range("O29").Value = 1.615
$ws.Range("P29").Value = 3.5
$ws.Range("Q29").Value = 4.5
$ws.Range("R29").Value = -0.75
$ws.Range("S29").Value = 1.85
$ws.Range("T29").Value = 1.95
$ws.Range("U29").Value = 2.75
$ws.Range("V29").Value = 2
$ws.Range("W29").Value = 1.8
$ws.Range("X29").Value = 0.615
$ws.Range("Y29").Value = -1
$ws.Range("Z29").Value = -1
$ws.Range("AA29").Value = 0.8500000000000001
$ws.Range("AB29").Value = -1
$ws.Range("AC29").Value = 1
$ws.Range("AD29").Value = -1

# Row 30
$ws.Range("B30").Value = 6865296
$ws.Range("C30").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D30").Value = 45172.61458333334
$ws.Range("E30").Value = 'Velez Mostar'
$ws.Range("F30").Value = 'Zeljeznicar'
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 'H'
$ws.Range("L30").Value = 1.909
$ws.Range("M30").Value = 3.2
$ws.Range("N30").Value = 3.6
$ws.Range("O30").Value = 1.95
$ws.Range("P30").Value = 3.2
$ws.Range("Q30").Value = 3.4
$ws.Range("R30").Value = -0.5
$ws.Range("S30").Value = 2.025
$ws.Range("T30").Value = 1.775
$ws.Range("U30").Value = 2.25
$ws.Range("V30").Value = 1.9
$ws.Range("W30").Value = 1.9
$ws.Range("X30").Value = 0.95
$ws.Range("Y30").Value = -1
$ws.Range("Z30").Value = -1
$ws.Range("AA30").Value = 1.025
$ws.Range("AB30").Value = -1
$ws.Range("AC30").Value = -1
$ws.Range("AD30").Value = 0.8999999999999999

# Row 111
$ws.Range("B111").Value = 6865352
$ws.Range("C111").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D111").Value = 45339.375
$ws.Range("E111").Value = 'NK Posusje'
$ws.Range("F111").Value = 'Zvijezda 09'
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 'H'
$ws.Range("L111").Value = 1.615
$ws.Range("M111").Value = 3.5
$ws.Range("N111").Value = 4.75
$ws.Range("O111").Value = 1.5
$ws.Range("P111").Value = 3.6
$ws.Range("Q111").Value = 5.75
$ws.Range("R111").Value = -1
$ws.Range("S111").Value = 1.9
$ws.Range("T111").Value = 1.9
$ws.Range("U111").Value = 2.25
$ws.Range("V111").Value = 1.85
$ws.Range("W111").Value = 1.95
$ws.Range("X111").Value = 0.5
$ws.Range("Y111").Value = -1
$ws.Range("Z111").Value = -1
$ws.Range("AA111").Value = 0.8999999999999999
$ws.Range("AB111").Value = -1
$ws.Range("AC111").Value = -0.5
$ws.Range("AD111").Value = 0.475

# Row 112
$ws.Range("B112").Value = 6865354
$ws.Range("C112").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D112").Value = 45339.375
$ws.Range("E112").Value = 'NK Igman Konjic'
$ws.Range("F112").Value = 'GOSK Gabela'
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 2
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1
$ws.Range("K112").Value = 'A'
$ws.Range("L112").Value = 1.8
$ws.Range("M112").Value = 3.25
$ws.Range("N112").Value = 4
$ws.Range("O112").Value = 2.25
$ws.Range("P112").Value = 3.1
$ws.Range("Q112").Value = 2.9
$ws.Range("R112").Value = -0.25
$ws.Range("S112").Value = 1.975
$ws.Range("T112").Value = 1.825
$ws.Range("U112").Value = 2.25
$ws.Range("V112").Value = 1.875
$ws.Range("W112").Value = 1.925
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = -1
$ws.Range("Z112").Value = 1.9
$ws.Range("AA112").Value = -1
$ws.Range("AB112").Value = 0.825
$ws.Range("AC112").Value = 0.875
$ws.Range("AD112").Value = -1

# Row 189
$ws.Range("B189").Value = 7952778
$ws.Range("C189").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D189").Value = 45432.5
$ws.Range("E189").Value = 'Sloga'
$ws.Range("F189").Value = 'Siroki Brijeg'
$ws.Range("G189").Value = 2
$ws.Range("H189").Value = 3
$ws.Range("I189").Value = 2
$ws.Range("J189").Value = 2
$ws.Range("K189").Value = 'A'
$ws.Range("L189").Value = 1.727
$ws.Range("M189").Value = 3.75
$ws.Range("N189").Value = 3.75
$ws.Range("O189").Value = 1.7
$ws.Range("P189").Value = 3.9
$ws.Range("Q189").Value = 3.9
$ws.Range("R189").Value = -0.75
$ws.Range("S189").Value = 1.975
$ws.Range("T189").Value = 1.825
$ws.Range("U189").Value = 2.25
$ws.Range("V189").Value = 1.8
$ws.Range("W189").Value = 2
$ws.Range("X189").Value = -1
$ws.Range("Y189").Value = -1
$ws.Range("Z189").Value = 2.9
$ws.Range("AA189").Value = -1
$ws.Range("AB189").Value = 0.825
$ws.Range("AC189").Value = 0.8
$ws.Range("AD189").Value = -1

# Row 190
$ws.Range("B190").Value = 7952780
$ws.Range("C190").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D190").Value = 45432.5
$ws.Range("E190").Value = 'Velez Mostar'
$ws.Range("F190").Value = 'GOSK Gabela'
$ws.Range("G190").Value = 3
$ws.Range("H190").Value = 3
$ws.Range("I190").Value = 1
$ws.Range("J190").Value = 1
$ws.Range("K190").Value = 'D'
$ws.Range("L190").Value = 1.4
$ws.Range("M190").Value = 4
$ws.Range("N190").Value = 7
$ws.Range("O190").Value = 1.363
$ws.Range("P190").Value = 4.2
$ws.Range("Q190").Value = 8
$ws.Range("R190").Value = -1.5
$ws.Range("S190").Value = 2
$ws.Range("T190").Value = 1.8
$ws.Range("U190").Value = 2.75
$ws.Range("V190").Value = 1.825
$ws.Range("W190").Value = 1.975
$ws.Range("X190").Value = -1
$ws.Range("Y190").Value = 3.2
$ws.Range("Z190").Value = -1
$ws.Range("AA190").Value = -1
$ws.Range("AB190").Value = 0.8
$ws.Range("AC190").Value = 0.825
$ws.Range("AD190").Value = -1

# Row 191
$ws.Range("B191").Value = 7952777
$ws.Range("C191").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D191").Value = 45432.5
$ws.Range("E191").Value = 'Borac Banja Luka'
$ws.Range("F191").Value = 'NK Igman Konjic'
$ws.Range("G191").Value = 4
$ws.Range("H191").Value = 3
$ws.Range("I191").Value = 1
$ws.Range("J191").Value = 2
$ws.Range("K191").Value = 'H'
$ws.Range("L191").Value = 1.25
$ws.Range("M191").Value = 5.75
$ws.Range("N191").Value = 7
$ws.Range("O191").Value = 1.2
$ws.Range("P191").Value = 5.75
$ws.Range("Q191").Value = 12
$ws.Range("R191").Value = -2
$ws.Range("S191").Value = 1.95
$ws.Range("T191").Value = 1.85
$ws.Range("U191").Value = 3.25
$ws.Range("V191").Value = 1.9
$ws.Range("W191").Value = 1.9
$ws.Range("X191").Value = 0.2
$ws.Range("Y191").Value = -1
$ws.Range("Z191").Value = -1
$ws.Range("AA191").Value = -1
$ws.Range("AB191").Value = 0.8500000000000001
$ws.Range("AC191").Value = 0.8999999999999999
$ws.Range("AD191").Value = -1

# Row 192
$ws.Range("B192").Value = 7952776
$ws.Range("C192").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D192").Value = 45432.5
$ws.Range("E192").Value = 'FK Sarajevo'
$ws.Range("F192").Value = 'NK Posusje'
$ws.Range("G192").Value = 1
$ws.Range("H192").Value = 1
$ws.Range("I192").Value = 0
$ws.Range("J192").Value = 0
$ws.Range("K192").Value = 'D'
$ws.Range("L192").Value = 1.571
$ws.Range("M192").Value = 3.4
$ws.Range("N192").Value = 5.5
$ws.Range("O192").Value = 1.363
$ws.Range("P192").Value = 3.9
$ws.Range("Q192").Value = 8
$ws.Range("R192").Value = -1.25
$ws.Range("S192").Value = 1.85
$ws.Range("T192").Value = 1.95
$ws.Range("U192").Value = 2.75
$ws.Range("V192").Value = 1.925
$ws.Range("W192").Value = 1.875
$ws.Range("X192").Value = -1
$ws.Range("Y192").Value = 2.9
$ws.Range("Z192").Value = -1
$ws.Range("AA192").Value = -1
$ws.Range("AB192").Value = 0.95
$ws.Range("AC192").Value = -1
$ws.Range("AD192").Value = 0.875

# Row 197
$ws.Range("B197").Value = 8259815
$ws.Range("C197").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D197").Value = 45438.5
$ws.Range("E197").Value = 'NK Posusje'
$ws.Range("F197").Value = 'Zvijezda 09'
$ws.Range("G197").Value = 2
$ws.Range("H197").Value = 0
$ws.Range("K197").Value = 'H'
$ws.Range("L197").Value = 1.4
$ws.Range("M197").Value = 4
$ws.Range("N197").Value = 6.5
$ws.Range("O197").Value = 1.25
$ws.Range("P197").Value = 5
$ws.Range("Q197").Value = 8.5
$ws.Range("R197").Value = -1.75
$ws.Range("S197").Value = 2
$ws.Range("T197").Value = 1.8
$ws.Range("U197").Value = 3
$ws.Range("V197").Value = 1.925
$ws.Range("W197").Value = 1.875
$ws.Range("X197").Value = 0.25
$ws.Range("Y197").Value = -1
$ws.Range("Z197").Value = -1
$ws.Range("AA197").Value = 0.5
$ws.Range("AB197").Value = -0.5
$ws.Range("AC197").Value = -1
$ws.Range("AD197").Value = 0.875

# Row 198
$ws.Range("B198").Value = 7952787
$ws.Range("C198").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D198").Value = 45438.5
$ws.Range("E198").Value = 'Zeljeznicar'
$ws.Range("F198").Value = 'Borac Banja Luka'
$ws.Range("G198").Value = 2
$ws.Range("H198").Value = 1
$ws.Range("K198").Value = 'H'
$ws.Range("L198").Value = 3.3
$ws.Range("M198").Value = 3
$ws.Range("N198").Value = 2.05
$ws.Range("O198").Value = 1.727
$ws.Range("P198").Value = 3.1
$ws.Range("Q198").Value = 4.5
$ws.Range("R198").Value = -0.5
$ws.Range("S198").Value = 1.825
$ws.Range("T198").Value = 1.975
$ws.Range("U198").Value = 2.25
$ws.Range("V198").Value = 2
$ws.Range("W198").Value = 1.8
$ws.Range("X198").Value = 0.7270000000000001
$ws.Range("Y198").Value = -1
$ws.Range("Z198").Value = -1
$ws.Range("AA198").Value = 0.825
$ws.Range("AB198").Value = -1
$ws.Range("AC198").Value = 1
$ws.Range("AD198").Value = -1

# Row 199
$ws.Range("B199").Value = 8259814
$ws.Range("C199").Value = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D199").Value = 45438.5
$ws.Range("E199").Value = 'Siroki Brijeg'
$ws.Range("F199").Value = 'FK Sarajevo'
$ws.Range("G199").Value = 2
$ws.Range("H199").Value = 2
$ws.Range("K199").Value = 'D'
$ws.Range("L199").Value = 3.4
$ws.Range("M199").Value = 3.1
$ws.Range("N199").Value = 2
$ws.Range("O199").Value = 9.5
$ws.Range("P199").Value = 4.75
$ws.Range("Q199").Value = 1.25
$ws.Range("R199").Value = 1.5
$ws.Range("S199").Value = 1.975
$ws.Range("T199").Value = 1.825
$ws.Range("U199").Value = 2.75
$ws.Range("V199").Value = 1.95
$ws.Range("W199").Value = 1.85
$ws.Range("X199").Value = -1
$ws.Range("Y199").Value = 3.75
$ws.Range("Z199").Value = -1
$ws.Range("AA199").Value = 0.9750000000000001
$ws.Range("AB199").Value = -1
$ws.Range("AC199").Value = 0.95
$ws.Range("AD199").Value = -1
